# 睡眠日记 Sleep Diary - add new week block (rows 158-175) and fill in
# day 6/7 columns (G/H) for the preceding week block (rows 142-155).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Duplicate the previous week's template block (rows 138-155, which
#    includes the "instructions" row, the name row, the day-label row,
#    the date row and the 14 question rows) down to rows 158-175. This
#    preserves styles, number formats and merged-cell geometry exactly
#    like the original author's copy/paste did.
# ---------------------------------------------------------------------
$ws.Range("A138:H155").Copy($ws.Range("A158"))

# Row heights are a row-level property and are not carried over by
# Range.Copy, so restore them to match the source rows.
$ws.Rows.Item(158).RowHeight = 17.25
$ws.Rows.Item(159).RowHeight = 16.5
$ws.Rows.Item(160).RowHeight = 16.5
$ws.Rows.Item(161).RowHeight = 18
$ws.Rows.Item(162).RowHeight = 33
$ws.Rows.Item(163).RowHeight = 16.5
$ws.Rows.Item(164).RowHeight = 16.5
$ws.Rows.Item(165).RowHeight = 16.5
$ws.Rows.Item(166).RowHeight = 49.5
$ws.Rows.Item(167).RowHeight = 16.5
$ws.Rows.Item(168).RowHeight = 33
$ws.Rows.Item(169).RowHeight = 33
$ws.Rows.Item(170).RowHeight = 82.5
$ws.Rows.Item(171).RowHeight = 66
$ws.Rows.Item(172).RowHeight = 49.5
$ws.Rows.Item(173).RowHeight = 99
$ws.Rows.Item(174).RowHeight = 115.5
$ws.Rows.Item(175).RowHeight = 66

# ---------------------------------------------------------------------
# 2) New week header: name stays "Minrui Ren" (already copied), the
#    7 day columns keep "第一天".."第七天" (already copied) and the date
#    row needs a fresh starting date (2026-01-17) plus "+1" formulas.
# ---------------------------------------------------------------------
$ws.Range("B161").Value = 46039
$ws.Range("C161").Formula = "=B161+1"
$ws.Range("D161").Formula = "=C161+1"
$ws.Range("E161").Formula = "=D161+1"
$ws.Range("F161").Formula = "=E161+1"
$ws.Range("G161").Formula = "=F161+1"
$ws.Range("H161").Formula = "=G161+1"

# ---------------------------------------------------------------------
# 3) Fill in the answers actually recorded for the new week: only the
#    first three days (columns B, C, D) were filled in by the author;
#    columns E-H stay blank (already cleared by the copy below).
# ---------------------------------------------------------------------

# Clear days 4-7 (E:H) for the new question rows - the source rows had
# data there (days 4 & 5) that must not bleed into the new week.
$ws.Range("E162:H175").ClearContents()

# 您今天早上几点醒来?
$ws.Range("B162").Value = "6：40"
$ws.Range("C162").Value = "6：30"
$ws.Range("D162").Value = "6：40"

# 您今天几点起床?
$ws.Range("B163").Value = "6：40"
$ws.Range("C163").Value = "6：30"
$ws.Range("D163").Value = "7：00"

# 您昨晚几点上床?
$ws.Range("B164").Value = "20：00"
$ws.Range("C164").Value = "23：20"
$ws.Range("D164").Value = "23：00"

# 您昨晚几点熄灯?
$ws.Range("B165").Value = "23：00"
$ws.Range("C165").Value = "23：20"
$ws.Range("D165").Value = "23：00"

# 您昨晚熄灯后花了多长时间入睡(分钟)?
$ws.Range("B166").Value = 0
$ws.Range("C166").Value = 10
$ws.Range("D166").Value = "10 min"

# 您整晚醒来几次?
$ws.Range("B167").Value = 1
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 1

# 您整晚总共醒了多长时间(分钟)?
$ws.Range("B168").Value = 30
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 5

# 您整晚总共睡了多长时间(分钟)?
$ws.Range("B169").Value = 430
$ws.Range("C169").Value = 420
$ws.Range("D169").Value = 460

# 您昨晚睡前是否使用了影响睡眠的物质...
$ws.Range("B170").Value = "无"
$ws.Range("C170").Value = "无"
$ws.Range("D170").Value = "无"

# 您昨晚睡前是否使用了电子产品?使用了多长时间(分钟)?
$ws.Range("B171").Value = 30
$ws.Range("C171").Value = 10
$ws.Range("D171").Value = 40

# 您昨晚睡前的身体紧张程度如何?
$ws.Range("B172").Value = 2
$ws.Range("C172").Value = 3
$ws.Range("D172").Value = 4

# 您昨晚睡前的精神紧张程度如何?
$ws.Range("B173").Value = 3
$ws.Range("C173").Value = 2
$ws.Range("D173").Value = 4

# 您整晚的睡眠质量如何?
$ws.Range("B174").Value = 3
$ws.Range("C174").Value = 2
$ws.Range("D174").Value = 4

# 您昨天白天是否小睡?
$ws.Range("B175").Value = "无"
$ws.Range("C175").Value = "无"
$ws.Range("D175").Value = "无"

# ---------------------------------------------------------------------
# 4) Back-fill the previous week (rows 142-155) day 6 & day 7 columns
#    (G & H), which had been left blank until now.
# ---------------------------------------------------------------------
$ws.Range("G142").Value = "6：20"
$ws.Range("H142").Value = "6：40"
$ws.Range("G143").Value = "6：20"
$ws.Range("H143").Value = "6：40"
$ws.Range("G144").Value = "24：00"
$ws.Range("H144").Value = "23：00"
$ws.Range("G145").Value = "24：00"
$ws.Range("H145").Value = "23：00"
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 40
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 1
$ws.Range("G148").Value = 5
$ws.Range("H148").Value = 5
$ws.Range("G149").Value = 380
$ws.Range("H149").Value = 420
$ws.Range("G150").Value = "无"
$ws.Range("H150").Value = "无"
$ws.Range("G151").Value = 20
$ws.Range("H151").Value = 40
$ws.Range("G152").Value = 4
$ws.Range("H152").Value = 4
$ws.Range("G153").Value = 4
$ws.Range("H153").Value = 4
$ws.Range("G154").Value = 4
$ws.Range("H154").Value = 4
$ws.Range("G155").Value = "无"
$ws.Range("H155").Value = "无"

# ---------------------------------------------------------------------
# 5) Leave the view scrolled/selected roughly where the author left it.
# ---------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A167"))
$ws.Range("C169").Select()

Write-Output "done"
